$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 15.74808933333333
$ws.Cells.Item(2, 8).Value = 47.244268
$ws.Cells.Item(2, 9).Value = 0.2237852983702856
$ws.Cells.Item(2, 10).Value = 0.2380674495116221
$ws.Cells.Item(2, 13).Value = 6.712486666666666
$ws.Cells.Item(2, 14).Value = 20.13746
$ws.Cells.Item(2, 15).Value = 0.6330487633990675
$ws.Cells.Item(2, 16).Value = 0.6414503882251803
$ws.Cells.Item(2, 17).Value = 105.7088396754755
$ws.Cells.Item(2, 18).Value = 951.3795570792798
$ws.Cells.Item(2, 19).Value = 0.1416670064002006
$ws.Cells.Item(2, 20).Value = 0.1527084579130085

$ws.Cells.Item(3, 7).Value = 15.74808933333333
$ws.Cells.Item(3, 8).Value = 47.244268
$ws.Cells.Item(3, 9).Value = 0.2237852983702856
$ws.Cells.Item(3, 10).Value = 0.2380674495116221
$ws.Cells.Item(3, 15).Value = 0.290741083484562
$ws.Cells.Item(3, 16).Value = 0.2945997080427384
$ws.Cells.Item(3, 17).Value = 48.54902869744534
$ws.Cells.Item(3, 18).Value = 436.941258277008
$ws.Cells.Item(3, 19).Value = 0.06506358011609281
$ws.Cells.Item(3, 20).Value = 0.07013460112060323

$ws.Cells.Item(4, 7).Value = 15.74808933333333
$ws.Cells.Item(4, 8).Value = 47.244268
$ws.Cells.Item(4, 9).Value = 0.2237852983702856
$ws.Cells.Item(4, 10).Value = 0.2380674495116221
$ws.Cells.Item(4, 13).Value = 0.2495096666666667
$ws.Cells.Item(4, 14).Value = 0.748529
$ws.Cells.Item(4, 15).Value = 0.02353103905946135
$ws.Cells.Item(4, 16).Value = 0.02384333563656022
$ws.Cells.Item(4, 17).Value = 3.929300520196889
$ws.Cells.Item(4, 18).Value = 35.363704681772
$ws.Cells.Item(4, 19).Value = 0.005265900596884403
$ws.Cells.Item(4, 20).Value = 0.005676322102845461

$ws.Cells.Item(5, 7).Value = 15.74808933333333
$ws.Cells.Item(5, 8).Value = 47.244268
$ws.Cells.Item(5, 9).Value = 0.2237852983702856
$ws.Cells.Item(5, 10).Value = 0.2380674495116221
$ws.Cells.Item(5, 13).Value = 0.4166465
$ws.Cells.Item(5, 14).Value = 0.8332930000000001
$ws.Cells.Item(5, 15).Value = 0.03929356804674715
$ws.Cells.Item(5, 16).Value = 0.02654337331298611
$ws.Cells.Item(5, 17).Value = 6.561386302420667
$ws.Cells.Item(5, 18).Value = 39.368317814524
$ws.Cells.Item(5, 19).Value = 0.008793322849374431
$ws.Cells.Item(5, 20).Value = 0.006319113186057458

$ws.Cells.Item(6, 7).Value = 15.74808933333333
$ws.Cells.Item(6, 8).Value = 47.244268
$ws.Cells.Item(6, 9).Value = 0.2237852983702856
$ws.Cells.Item(6, 10).Value = 0.2380674495116221
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.1419326666666667
$ws.Cells.Item(6, 14).Value = 0.425798
$ws.Cells.Item(6, 15).Value = 0.01338554601016197
$ws.Cells.Item(6, 16).Value = 0.01356319478253491
$ws.Cells.Item(6, 17).Value = 2.235168313984889
$ws.Cells.Item(6, 18).Value = 20.116514825864
$ws.Cells.Item(6, 19).Value = 0.002995488407733282
$ws.Cells.Item(6, 20).Value = 0.003228955189107425

$ws.Cells.Item(7, 9).Value = 0.366480229475165
$ws.Cells.Item(7, 10).Value = 0.38986928168634
$ws.Cells.Item(7, 13).Value = 6.712486666666666
$ws.Cells.Item(7, 14).Value = 20.13746
$ws.Cells.Item(7, 15).Value = 0.6330487633990675
$ws.Cells.Item(7, 16).Value = 0.6414503882251803
$ws.Cells.Item(7, 17).Value = 173.1132478493755
$ws.Cells.Item(7, 18).Value = 1558.01923064438
$ws.Cells.Item(7, 19).Value = 0.2319998560794597
$ws.Cells.Item(7, 20).Value = 0.250081802094775

$ws.Cells.Item(8, 9).Value = 0.366480229475165
$ws.Cells.Item(8, 10).Value = 0.38986928168634
$ws.Cells.Item(8, 15).Value = 0.290741083484562
$ws.Cells.Item(8, 16).Value = 0.2945997080427384
$ws.Cells.Item(8, 19).Value = 0.1065508589932804
$ws.Cells.Item(8, 20).Value = 0.1148553765596279

$ws.Cells.Item(9, 9).Value = 0.366480229475165
$ws.Cells.Item(9, 10).Value = 0.38986928168634
$ws.Cells.Item(9, 13).Value = 0.2495096666666667
$ws.Cells.Item(9, 14).Value = 0.748529
$ws.Cells.Item(9, 15).Value = 0.02353103905946135
$ws.Cells.Item(9, 16).Value = 0.02384333563656022
$ws.Cells.Item(9, 17).Value = 6.434788016931888
$ws.Cells.Item(9, 18).Value = 57.913092152387
$ws.Cells.Item(9, 19).Value = 0.008623660594300468
$ws.Cells.Item(9, 20).Value = 0.009295784137632048

$ws.Cells.Item(10, 9).Value = 0.366480229475165
$ws.Cells.Item(10, 10).Value = 0.38986928168634
$ws.Cells.Item(10, 13).Value = 0.4166465
$ws.Cells.Item(10, 14).Value = 0.8332930000000001
$ws.Cells.Item(10, 15).Value = 0.03929356804674715
$ws.Cells.Item(10, 16).Value = 0.02654337331298611
$ws.Cells.Item(10, 17).Value = 10.74520254591317
$ws.Cells.Item(10, 18).Value = 64.471215275479
$ws.Cells.Item(10, 19).Value = 0.01440031583466991
$ws.Cells.Item(10, 20).Value = 0.01034844588706626

$ws.Cells.Item(11, 9).Value = 0.366480229475165
$ws.Cells.Item(11, 10).Value = 0.38986928168634
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.1419326666666667
$ws.Cells.Item(11, 14).Value = 0.425798
$ws.Cells.Item(11, 15).Value = 0.01338554601016197
$ws.Cells.Item(11, 16).Value = 0.01356319478253491
$ws.Cells.Item(11, 17).Value = 3.660405766554889
$ws.Cells.Item(11, 18).Value = 32.943651898994
$ws.Cells.Item(11, 19).Value = 0.004905537973454537
$ws.Cells.Item(11, 20).Value = 0.005287873007238798

$ws.Cells.Item(12, 7).Value = 2.416095
$ws.Cells.Item(12, 8).Value = 7.248285
$ws.Cells.Item(12, 9).Value = 0.03433346922420018
$ws.Cells.Item(12, 10).Value = 0.03652465783327086
$ws.Cells.Item(12, 13).Value = 6.712486666666666
$ws.Cells.Item(12, 14).Value = 20.13746
$ws.Cells.Item(12, 15).Value = 0.6330487633990675
$ws.Cells.Item(12, 16).Value = 0.6414503882251803
$ws.Cells.Item(12, 17).Value = 16.2180054729
$ws.Cells.Item(12, 18).Value = 145.9620492561
$ws.Cells.Item(12, 19).Value = 0.02173476023557987
$ws.Cells.Item(12, 20).Value = 0.02342875594694346

$ws.Cells.Item(13, 7).Value = 2.416095
$ws.Cells.Item(13, 8).Value = 7.248285
$ws.Cells.Item(13, 9).Value = 0.03433346922420018
$ws.Cells.Item(13, 10).Value = 0.03652465783327086
$ws.Cells.Item(13, 15).Value = 0.290741083484562
$ws.Cells.Item(13, 16).Value = 0.2945997080427384
$ws.Cells.Item(13, 17).Value = 7.44846330294
$ws.Cells.Item(13, 18).Value = 67.03616972646
$ws.Cells.Item(13, 19).Value = 0.009982150042027824
$ws.Cells.Item(13, 20).Value = 0.01076015353404251

$ws.Cells.Item(14, 7).Value = 2.416095
$ws.Cells.Item(14, 8).Value = 7.248285
$ws.Cells.Item(14, 9).Value = 0.03433346922420018
$ws.Cells.Item(14, 10).Value = 0.03652465783327086
$ws.Cells.Item(14, 13).Value = 0.2495096666666667
$ws.Cells.Item(14, 14).Value = 0.748529
$ws.Cells.Item(14, 15).Value = 0.02353103905946135
$ws.Cells.Item(14, 16).Value = 0.02384333563656022
$ws.Cells.Item(14, 17).Value = 0.6028390580849999
$ws.Cells.Item(14, 18).Value = 5.425551522765
$ws.Cells.Item(14, 19).Value = 0.0008079022053614688
$ws.Cells.Item(14, 20).Value = 0.0008708696757291956

$ws.Cells.Item(15, 7).Value = 2.416095
$ws.Cells.Item(15, 8).Value = 7.248285
$ws.Cells.Item(15, 9).Value = 0.03433346922420018
$ws.Cells.Item(15, 10).Value = 0.03652465783327086
$ws.Cells.Item(15, 13).Value = 0.4166465
$ws.Cells.Item(15, 14).Value = 0.8332930000000001
$ws.Cells.Item(15, 15).Value = 0.03929356804674715
$ws.Cells.Item(15, 16).Value = 0.02654337331298611
$ws.Cells.Item(15, 17).Value = 1.0066575254175
$ws.Cells.Item(15, 18).Value = 6.039945152505
$ws.Cells.Item(15, 19).Value = 0.001349084509242009
$ws.Cells.Item(15, 20).Value = 0.0009694876279975908

$ws.Cells.Item(16, 7).Value = 2.416095
$ws.Cells.Item(16, 8).Value = 7.248285
$ws.Cells.Item(16, 9).Value = 0.03433346922420018
$ws.Cells.Item(16, 10).Value = 0.03652465783327086
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.1419326666666667
$ws.Cells.Item(16, 14).Value = 0.425798
$ws.Cells.Item(16, 15).Value = 0.01338554601016197
$ws.Cells.Item(16, 16).Value = 0.01356319478253491
$ws.Cells.Item(16, 17).Value = 0.34292280627
$ws.Cells.Item(16, 18).Value = 3.08630525643
$ws.Cells.Item(16, 19).Value = 0.0004595722319890114
$ws.Cells.Item(16, 20).Value = 0.000495391048558092

$ws.Cells.Item(17, 7).Value = 12.6651745
$ws.Cells.Item(17, 8).Value = 25.330349
$ws.Cells.Item(17, 9).Value = 0.1799761097617747
$ws.Cells.Item(17, 10).Value = 0.1276415496938013
$ws.Cells.Item(17, 13).Value = 6.712486666666666
$ws.Cells.Item(17, 14).Value = 20.13746
$ws.Cells.Item(17, 15).Value = 0.6330487633990675
$ws.Cells.Item(17, 16).Value = 0.6414503882251803
$ws.Cells.Item(17, 17).Value = 85.01481496225665
$ws.Cells.Item(17, 18).Value = 510.0888897735399
$ws.Cells.Item(17, 19).Value = 0.1139336537260663
$ws.Cells.Item(17, 20).Value = 0.08187572160475248

$ws.Cells.Item(18, 7).Value = 12.6651745
$ws.Cells.Item(18, 8).Value = 25.330349
$ws.Cells.Item(18, 9).Value = 0.1799761097617747
$ws.Cells.Item(18, 10).Value = 0.1276415496938013
$ws.Cells.Item(18, 15).Value = 0.290741083484562
$ws.Cells.Item(18, 16).Value = 0.2945997080427384
$ws.Cells.Item(18, 17).Value = 39.044858537674
$ws.Cells.Item(18, 18).Value = 234.269151226044
$ws.Cells.Item(18, 19).Value = 0.05232644915347481
$ws.Cells.Item(18, 20).Value = 0.03760316327391654

$ws.Cells.Item(19, 7).Value = 12.6651745
$ws.Cells.Item(19, 8).Value = 25.330349
$ws.Cells.Item(19, 9).Value = 0.1799761097617747
$ws.Cells.Item(19, 10).Value = 0.1276415496938013
$ws.Cells.Item(19, 13).Value = 0.2495096666666667
$ws.Cells.Item(19, 14).Value = 0.748529
$ws.Cells.Item(19, 15).Value = 0.02353103905946135
$ws.Cells.Item(19, 16).Value = 0.02384333563656022
$ws.Cells.Item(19, 17).Value = 3.160083467770166
$ws.Cells.Item(19, 18).Value = 18.960500806621
$ws.Cells.Item(19, 19).Value = 0.004235024868574223
$ws.Cells.Item(19, 20).Value = 0.003043400310519985

$ws.Cells.Item(20, 7).Value = 12.6651745
$ws.Cells.Item(20, 8).Value = 25.330349
$ws.Cells.Item(20, 9).Value = 0.1799761097617747
$ws.Cells.Item(20, 10).Value = 0.1276415496938013
$ws.Cells.Item(20, 13).Value = 0.4166465
$ws.Cells.Item(20, 14).Value = 0.8332930000000001
$ws.Cells.Item(20, 15).Value = 0.03929356804674715
$ws.Cells.Item(20, 16).Value = 0.02654337331298611
$ws.Cells.Item(20, 17).Value = 5.27690062731425
$ws.Cells.Item(20, 18).Value = 21.107602509257
$ws.Cells.Item(20, 19).Value = 0.007071903515713127
$ws.Cells.Item(20, 20).Value = 0.003388037303770636

$ws.Cells.Item(21, 7).Value = 12.6651745
$ws.Cells.Item(21, 8).Value = 25.330349
$ws.Cells.Item(21, 9).Value = 0.1799761097617747
$ws.Cells.Item(21, 10).Value = 0.1276415496938013
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 12).Value = 0.6666666666666666
$ws.Cells.Item(21, 13).Value = 0.1419326666666667
$ws.Cells.Item(21, 14).Value = 0.425798
$ws.Cells.Item(21, 15).Value = 0.01338554601016197
$ws.Cells.Item(21, 16).Value = 0.01356319478253491
$ws.Cells.Item(21, 17).Value = 1.797601990583667
$ws.Cells.Item(21, 18).Value = 10.785611943502
$ws.Cells.Item(21, 19).Value = 0.002409078497946195
$ws.Cells.Item(21, 20).Value = 0.001731227200841636

$ws.Cells.Item(22, 7).Value = 13.75232733333333
$ws.Cells.Item(22, 8).Value = 41.25698199999999
$ws.Cells.Item(22, 9).Value = 0.1954248931685745
$ws.Cells.Item(22, 10).Value = 0.2078970612749657
$ws.Cells.Item(22, 13).Value = 6.712486666666666
$ws.Cells.Item(22, 14).Value = 20.13746
$ws.Cells.Item(22, 15).Value = 0.6330487633990675
$ws.Cells.Item(22, 16).Value = 0.6414503882251803
$ws.Cells.Item(22, 17).Value = 92.31231386063554
$ws.Cells.Item(22, 18).Value = 830.8108247457197
$ws.Cells.Item(22, 19).Value = 0.1237134869577609
$ws.Cells.Item(22, 20).Value = 0.1333556506657008

$ws.Cells.Item(23, 7).Value = 13.75232733333333
$ws.Cells.Item(23, 8).Value = 41.25698199999999
$ws.Cells.Item(23, 9).Value = 0.1954248931685745
$ws.Cells.Item(23, 10).Value = 0.2078970612749657
$ws.Cells.Item(23, 15).Value = 0.290741083484562
$ws.Cells.Item(23, 16).Value = 0.2945997080427384
$ws.Cells.Item(23, 17).Value = 42.39638982422133
$ws.Cells.Item(23, 18).Value = 381.567508417992
$ws.Cells.Item(23, 19).Value = 0.05681804517968611
$ws.Cells.Item(23, 20).Value = 0.06124641355454817

$ws.Cells.Item(24, 7).Value = 13.75232733333333
$ws.Cells.Item(24, 8).Value = 41.25698199999999
$ws.Cells.Item(24, 9).Value = 0.1954248931685745
$ws.Cells.Item(24, 10).Value = 0.2078970612749657
$ws.Cells.Item(24, 13).Value = 0.2495096666666667
$ws.Cells.Item(24, 14).Value = 0.748529
$ws.Cells.Item(24, 15).Value = 0.02353103905946135
$ws.Cells.Item(24, 16).Value = 0.02384333563656022
$ws.Cells.Item(24, 17).Value = 3.431338608830889
$ws.Cells.Item(24, 18).Value = 30.882047479478
$ws.Cells.Item(24, 19).Value = 0.004598550794340788
$ws.Cells.Item(24, 20).Value = 0.004956959409833533

$ws.Cells.Item(25, 7).Value = 13.75232733333333
$ws.Cells.Item(25, 8).Value = 41.25698199999999
$ws.Cells.Item(25, 9).Value = 0.1954248931685745
$ws.Cells.Item(25, 10).Value = 0.2078970612749657
$ws.Cells.Item(25, 13).Value = 0.4166465
$ws.Cells.Item(25, 14).Value = 0.8332930000000001
$ws.Cells.Item(25, 15).Value = 0.03929356804674715
$ws.Cells.Item(25, 16).Value = 0.02654337331298611
$ws.Cells.Item(25, 17).Value = 5.729859050287667
$ws.Cells.Item(25, 18).Value = 34.379154301726
$ws.Cells.Item(25, 19).Value = 0.007678941337747673
$ws.Cells.Item(25, 20).Value = 0.005518289308094161

$ws.Cells.Item(26, 7).Value = 13.75232733333333
$ws.Cells.Item(26, 8).Value = 41.25698199999999
$ws.Cells.Item(26, 9).Value = 0.1954248931685745
$ws.Cells.Item(26, 10).Value = 0.2078970612749657
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 12).Value = 0.6666666666666666
$ws.Cells.Item(26, 13).Value = 0.1419326666666667
$ws.Cells.Item(26, 14).Value = 0.425798
$ws.Cells.Item(26, 15).Value = 0.01338554601016197
$ws.Cells.Item(26, 16).Value = 0.01356319478253491
$ws.Cells.Item(26, 17).Value = 1.951904491292889
$ws.Cells.Item(26, 18).Value = 17.567140421636
$ws.Cells.Item(26, 19).Value = 0.00261586889903894
$ws.Cells.Item(26, 20).Value = 0.002819748336788954
